$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, matching the bold/bordered header style
# used by the other header cells (B1:E1) by copying E1's formatting onto F1.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Populate F2:F59 with the recorded time_taken timestamps (one per data row),
# using plain (unstyled) cells like the rest of the data columns.
$timeTaken = @(
    "2021-10-05 13:38:35.886479",
    "2021-10-05 13:38:35.886492",
    "2021-10-05 13:38:35.886497",
    "2021-10-05 13:38:35.886500",
    "2021-10-05 13:38:35.886503",
    "2021-10-05 13:38:35.886506",
    "2021-10-05 13:38:35.886509",
    "2021-10-05 13:38:35.886512",
    "2021-10-05 13:38:35.886516",
    "2021-10-05 13:38:35.886519",
    "2021-10-05 13:38:35.886522",
    "2021-10-05 13:38:35.886525",
    "2021-10-05 13:38:35.886527",
    "2021-10-05 13:38:35.886530",
    "2021-10-05 13:38:35.886533",
    "2021-10-05 13:38:35.886536",
    "2021-10-05 13:38:35.886539",
    "2021-10-05 13:38:35.886543",
    "2021-10-05 13:38:35.886545",
    "2021-10-05 13:38:35.886548",
    "2021-10-05 13:38:35.886551",
    "2021-10-05 13:38:35.886554",
    "2021-10-05 13:38:35.886557",
    "2021-10-05 13:38:35.886560",
    "2021-10-05 13:38:35.886563",
    "2021-10-05 13:38:35.886566",
    "2021-10-05 13:38:35.886569",
    "2021-10-05 13:38:35.886572",
    "2021-10-05 13:38:35.886575",
    "2021-10-05 13:38:35.886578",
    "2021-10-05 13:38:35.886581",
    "2021-10-05 13:38:35.886584",
    "2021-10-05 13:38:35.886587",
    "2021-10-05 13:38:35.886590",
    "2021-10-05 13:38:35.886593",
    "2021-10-05 13:38:35.886596",
    "2021-10-05 13:38:35.886599",
    "2021-10-05 13:38:35.886602",
    "2021-10-05 13:38:35.886605",
    "2021-10-05 13:38:35.886608",
    "2021-10-05 13:38:35.886611",
    "2021-10-05 13:38:35.886614",
    "2021-10-05 13:38:35.886617",
    "2021-10-05 13:38:35.886620",
    "2021-10-05 13:38:35.886623",
    "2021-10-05 13:38:35.886626",
    "2021-10-05 13:38:35.886629",
    "2021-10-05 13:38:35.886631",
    "2021-10-05 13:38:35.886634",
    "2021-10-05 13:38:35.886637",
    "2021-10-05 13:38:35.886640",
    "2021-10-05 13:38:35.886643",
    "2021-10-05 13:38:35.886646",
    "2021-10-05 13:38:35.886649",
    "2021-10-05 13:38:35.886652",
    "2021-10-05 13:38:35.886655",
    "2021-10-05 13:38:35.886658",
    "2021-10-05 13:38:35.886661"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $timeTaken[$i]
}
